$wb = $excel.ActiveWorkbook

# Sheet ALC, row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5004.3125
$ws.Range("I40").Value = 3292
$ws.Range("K40").Value = 3292
$ws.Range("M40").Value = -3117

# Sheet ALC, row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 5967.7144
$ws.Range("I43").Value = 4999
$ws.Range("J43").Value = 6694.25
$ws.Range("K43").Value = 4999
$ws.Range("L43").Value = 6694.25
$ws.Range("M43").Value = -4930
$ws.Range("N43").Value = -6832.25

# Sheet ALC, row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1990
$ws.Range("I62").Value = 1990
$ws.Range("K62").Value = 1990
$ws.Range("M62").Value = -1366

# Sheet ALC, row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 1990
$ws.Range("I65").Value = 1990
$ws.Range("K65").Value = 9950
$ws.Range("M65").Value = -6830

# Sheet ALC, row 80
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 2495.923
$ws.Range("J80").Value = 3144.4443
$ws.Range("L80").Value = 9433.332900000001
$ws.Range("N80").Value = -11429.3329

# Sheet ALC, row 83
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 2495.923
$ws.Range("J83").Value = 3144.4443
$ws.Range("L83").Value = 28299.9987
$ws.Range("N83").Value = -38283.9987

# Sheet ALC, row 88
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 2435.625
$ws.Range("J88").Value = 2319.7144
$ws.Range("L88").Value = 2319.7144
$ws.Range("N88").Value = -3131.7144

# Sheet ALC, row 91
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 2435.625
$ws.Range("J91").Value = 2319.7144
$ws.Range("L91").Value = 2319.7144
$ws.Range("N91").Value = -5127.7144

# Sheet ALC, row 97
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 2806
$ws.Range("J97").Value = 2806
$ws.Range("L97").Value = 8418
$ws.Range("N97").Value = -9410

# Sheet ALC, row 103
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 7497.5
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 7497.5
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 22492.5
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = -23664.5

# Sheet ALC, row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 3212.8
$ws.Range("J112").Value = 3212.8
$ws.Range("L112").Value = 9638.400000000001
$ws.Range("N112").Value = -11854.4

# Sheet ALC, row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 5185.2
$ws.Range("I116").Value = 5106.75
$ws.Range("K116").Value = 5106.75
$ws.Range("M116").Value = -1664.75

# Sheet ALC, row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 2428.077
$ws.Range("I125").Value = 2005.1666
$ws.Range("J125").Value = 2790.5715
$ws.Range("K125").Value = 18046.4994
$ws.Range("L125").Value = 25115.1435
$ws.Range("M125").Value = -15586.4994
$ws.Range("N125").Value = -30035.1435

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3312.652
$ws.Range("I138").Value = 1462.125
$ws.Range("K138").Value = 4386.375
$ws.Range("M138").Value = 753.625

# Sheet ARM, row 36
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

# Sheet ARM, row 41
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 2713.5
$ws.Range("I41").Value = 1951.3334
$ws.Range("J41").Value = 5000
$ws.Range("K41").Value = 1951.3334
$ws.Range("L41").Value = 5000
$ws.Range("M41").Value = -1537.3334
$ws.Range("N41").Value = -5828

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2053.375
$ws.Range("I74").Value = 1990.2667
$ws.Range("K74").Value = 1990.2667
$ws.Range("M74").Value = -1116.2667

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2053.375
$ws.Range("I77").Value = 1990.2667
$ws.Range("K77").Value = 9951.333499999999
$ws.Range("M77").Value = -5583.333499999999

# Sheet ARM, row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1810.1875
$ws.Range("I88").Value = 1359.125
$ws.Range("K88").Value = 1359.125
$ws.Range("M88").Value = -953.125

# Sheet ARM, row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1810.1875
$ws.Range("I91").Value = 1359.125
$ws.Range("K91").Value = 1359.125
$ws.Range("M91").Value = 44.875

# Sheet ARM, row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3225.75
$ws.Range("I132").Value = 3225.75
$ws.Range("K132").Value = 9677.25
$ws.Range("M132").Value = -7147.25

# Sheet BSM, row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1006.4
$ws.Range("I105").Value = 1007
$ws.Range("J105").Value = 1005.5
$ws.Range("K105").Value = 1007
$ws.Range("L105").Value = 1005.5
$ws.Range("M105").Value = 740
$ws.Range("N105").Value = -4499.5

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1339.6
$ws.Range("I134").Value = 1339.6
$ws.Range("K134").Value = 4018.8
$ws.Range("M134").Value = -1483.8

# Sheet CRP, row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3879.6
$ws.Range("I58").Value = 4301.5
$ws.Range("J58").Value = 2192
$ws.Range("K58").Value = 4301.5
$ws.Range("L58").Value = 2192
$ws.Range("M58").Value = -4098.5
$ws.Range("N58").Value = -2598

# Sheet CRP, row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 998002
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

# Sheet CRP, row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 998002
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

# Sheet CRP, row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1244.8572
$ws.Range("I107").Value = 928.34784
$ws.Range("K107").Value = 928.34784
$ws.Range("M107").Value = 991.65216

# Sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1798.4073
$ws.Range("I132").Value = 1941.2858
$ws.Range("J132").Value = 1298.3334
$ws.Range("K132").Value = 5823.857400000001
$ws.Range("L132").Value = 3895.0002
$ws.Range("M132").Value = -3293.857400000001
$ws.Range("N132").Value = -8955.0002

# Sheet CRP, row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4505
$ws.Range("I134").Value = 4505
$ws.Range("K134").Value = 13515
$ws.Range("M134").Value = -10980

# Sheet CRP, row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3879.6
$ws.Range("I136").Value = 4301.5
$ws.Range("J136").Value = 2192
$ws.Range("K136").Value = 12904.5
$ws.Range("L136").Value = 6576
$ws.Range("M136").Value = -10354.5
$ws.Range("N136").Value = -11676

# Sheet CUL, row 80
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4527
$ws.Range("J80").Value = 4866.3335
$ws.Range("L80").Value = 14599.0005
$ws.Range("N80").Value = -16471.0005

# Sheet CUL, row 83
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 4527
$ws.Range("J83").Value = 4866.3335
$ws.Range("L83").Value = 43797.0015
$ws.Range("N83").Value = -53157.0015

# Sheet CUL, row 93
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 550
$ws.Range("J93").Value = 900
$ws.Range("L93").Value = 2700
$ws.Range("N93").Value = -6444

# Sheet CUL, row 128
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 277418
$ws.Range("I128").Value = 277418
$ws.Range("K128").Value = 832254
$ws.Range("M128").Value = -827274

# Sheet CUL, row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3932.3333
$ws.Range("J137").Value = 3932.3333
$ws.Range("L137").Value = 11796.9999
$ws.Range("N137").Value = -21996.9999

# Sheet GSM, row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 443.0909
$ws.Range("I2").Value = 121.5
$ws.Range("J2").Value = 626.8570999999999
$ws.Range("K2").Value = 121.5
$ws.Range("L2").Value = 626.8570999999999
$ws.Range("M2").Value = -8.5
$ws.Range("N2").Value = -852.8570999999999

# Sheet GSM, row 46
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 13333.333
$ws.Range("J46").Value = 13333.333
$ws.Range("L46").Value = 13333.333
$ws.Range("N46").Value = -13645.333

# Sheet GSM, row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 7426.857
$ws.Range("I107").Value = 297.25
$ws.Range("J107").Value = 16933
$ws.Range("K107").Value = 297.25
$ws.Range("L107").Value = 16933
$ws.Range("M107").Value = 1622.75
$ws.Range("N107").Value = -20773

# Sheet LTW, row 32
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

# Sheet LTW, row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 705
$ws.Range("I61").Value = 705
$ws.Range("K61").Value = 705
$ws.Range("M61").Value = -503

# Sheet LTW, row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1225.6666
$ws.Range("I82").Value = 1240.8334
$ws.Range("K82").Value = 1240.8334
$ws.Range("M82").Value = -879.8334

# Sheet LTW, row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1225.6666
$ws.Range("I85").Value = 1240.8334
$ws.Range("K85").Value = 1240.8334
$ws.Range("M85").Value = 7.166600000000017

# Sheet LTW, row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 705
$ws.Range("I113").Value = 705
$ws.Range("K113").Value = 705
$ws.Range("M113").Value = 1465

# Sheet WVR, row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 625
$ws.Range("I107").Value = 160
$ws.Range("J107").Value = 1090
$ws.Range("K107").Value = 480
$ws.Range("L107").Value = 3270
$ws.Range("M107").Value = 1440
$ws.Range("N107").Value = -7110

# Sheet WVR, row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 461.89474
$ws.Range("I113").Value = 427.36365
$ws.Range("J113").Value = 509.375
$ws.Range("K113").Value = 1282.09095
$ws.Range("L113").Value = 1528.125
$ws.Range("M113").Value = 887.90905
$ws.Range("N113").Value = -5868.125

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8000
$ws.Range("I132").Value = 8000
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 24000
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -21470
$ws.Range("N132").ClearContents()

# Sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1000.26666
$ws.Range("I136").Value = 958.75
$ws.Range("K136").Value = 2876.25
$ws.Range("M136").Value = -326.25
